$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "notifications"

# Update Steps (column D) text for rows 4-16 with new granular steps
$ws.Cells.Item(4, 4).Value = "1. Open the app`n2. login to the app using premium users`n3. tap on the settings`n4. find the notification toogle button"
$ws.Cells.Item(5, 4).Value = "1. Open the app`n2. login to the app using premium users`n3. tap on the settings`n4. find and observe the notification toogle button"
$ws.Cells.Item(6, 4).Value = "1. Open the app`n2. login to the app using premium users`n3. tap on the settings`n4. find and observe the notification toogle button"
$ws.Cells.Item(7, 4).Value = "1. Open the app`n2. login to the app using premium users`n3. tap on the settings`n4. find and tap on the notification toogle button"
$ws.Cells.Item(8, 4).Value = "1. Open the app`n2. login to the app using premium users`n3. tap on the settings`n4. find and tap on the notification toogle button"
$ws.Cells.Item(9, 4).Value = "1. Open the app`n2. login to the app using premium users`n3. tap on the settings`n4. find and tap on the notification toogle button"
$ws.Cells.Item(10, 4).Value = "1. Open the app`n2. login to the app using premium users`n3. tap on the settings`n4. find and disable the notification toogle button"
$ws.Cells.Item(11, 4).Value = "1. Open the app`n2. login to the app using premium users`n3. tap on the settings`n4. find and disable the notification toogle button"
$ws.Cells.Item(12, 4).Value = "1. Open the app`n2. login to the app using premium users`n3. tap on the settings`n4. find and disable the notification toogle button`n5. close the app.`n6. open the app."
$ws.Cells.Item(13, 4).Value = "1. Open the app`n2. login to the app using premium users`n3. tap on the settings`n4. find and enable the notification toogle button"
$ws.Cells.Item(14, 4).Value = "1. Open the app`n2. login to the app using premium users`n3. tap on the settings`n4. find and disable the notification toogle button"
$ws.Cells.Item(15, 4).Value = "1. Open the app`n2. login to the app using premium users`n3. tap on the settings`n4. find and enable the notification toogle button`n5. try to enable again"
$ws.Cells.Item(16, 4).Value = "1. Open the app`n2. login to the app using premium users`n3. tap on the settings`n4. find and disable the notification toogle button`n5. try to disable again"

# Update Test Case ID (column A) for rows 31-57 to new unique sequential IDs
$ws.Cells.Item(31, 1).Value = "SYMENADISNOT-023"
$ws.Cells.Item(32, 1).Value = "SYMENADISNOT-024"
$ws.Cells.Item(33, 1).Value = "SYMENADISNOT-025"
$ws.Cells.Item(34, 1).Value = "SYMENADISNOT-026"
$ws.Cells.Item(35, 1).Value = "SYMENADISNOT-027"
$ws.Cells.Item(36, 1).Value = "SYMENADISNOT-028"
$ws.Cells.Item(37, 1).Value = "SYMENADISNOT-029"
$ws.Cells.Item(38, 1).Value = "SYMENADISNOT-030"
$ws.Cells.Item(39, 1).Value = "SYMENADISNOT-031"
$ws.Cells.Item(40, 1).Value = "SYMENADISNOT-032"
$ws.Cells.Item(41, 1).Value = "SYMENADISNOT-033"
$ws.Cells.Item(42, 1).Value = "SYMENADISNOT-034"
$ws.Cells.Item(43, 1).Value = "SYMENADISNOT-035"
$ws.Cells.Item(44, 1).Value = "SYMENADISNOT-036"
$ws.Cells.Item(45, 1).Value = "SYMENADISNOT-037"
$ws.Cells.Item(46, 1).Value = "SYMENADISNOT-038"
$ws.Cells.Item(47, 1).Value = "SYMENADISNOT-039"
$ws.Cells.Item(48, 1).Value = "SYMENADISNOT-040"
$ws.Cells.Item(49, 1).Value = "SYMENADISNOT-041"
$ws.Cells.Item(50, 1).Value = "SYMENADISNOT-042"
$ws.Cells.Item(51, 1).Value = "SYMENADISNOT-043"
$ws.Cells.Item(52, 1).Value = "SYMENADISNOT-044"
$ws.Cells.Item(53, 1).Value = "SYMENADISNOT-045"
$ws.Cells.Item(54, 1).Value = "SYMENADISNOT-046"
$ws.Cells.Item(55, 1).Value = "SYMENADISNOT-047"
$ws.Cells.Item(56, 1).Value = "SYMENADISNOT-048"
$ws.Cells.Item(57, 1).Value = "SYMENADISNOT-049"

# Row height adjustments (Excel auto-fit results captured from the authored edit)
$ws.Rows.Item(4).RowHeight = 58.8
$ws.Rows.Item(5).RowHeight = 64.2
$ws.Rows.Item(6).RowHeight = 64.8
$ws.Rows.Item(7).RowHeight = 61.2
$ws.Rows.Item(8).RowHeight = 68.4
$ws.Rows.Item(9).RowHeight = 73.2
$ws.Rows.Item(10).RowHeight = 66.6
$ws.Rows.Item(11).RowHeight = 67.8
$ws.Rows.Item(12).RowHeight = 88.8
$ws.Rows.Item(13).RowHeight = 61.8
$ws.Rows.Item(14).RowHeight = 70.2
$ws.Rows.Item(15).RowHeight = 78.6
$ws.Rows.Item(16).RowHeight = 75
$ws.Rows.Item(17).RowHeight = 52.2
$ws.Rows.Item(18).RowHeight = 60.6
$ws.Rows.Item(19).RowHeight = 51
$ws.Rows.Item(20).RowHeight = 52.2
$ws.Rows.Item(21).RowHeight = 29.4
$ws.Rows.Item(22).RowHeight = 45.6
$ws.Rows.Item(23).RowHeight = 55.2
$ws.Rows.Item(24).RowHeight = 52.2
$ws.Rows.Item(25).RowHeight = 55.2
$ws.Rows.Item(26).RowHeight = 60
$ws.Rows.Item(27).RowHeight = 50.4
$ws.Rows.Item(28).RowHeight = 53.4
$ws.Rows.Item(29).RowHeight = 37.2
$ws.Rows.Item(30).RowHeight = 46.8
$ws.Rows.Item(31).RowHeight = 51
$ws.Rows.Item(32).RowHeight = 51
$ws.Rows.Item(33).RowHeight = 45.6
$ws.Rows.Item(34).RowHeight = 49.8
$ws.Rows.Item(35).RowHeight = 59.4
$ws.Rows.Item(36).RowHeight = 48.6
$ws.Rows.Item(37).RowHeight = 55.2
$ws.Rows.Item(38).RowHeight = 48.6
$ws.Rows.Item(39).RowHeight = 49.8
$ws.Rows.Item(40).RowHeight = 45.6
$ws.Rows.Item(41).RowHeight = 43.8
$ws.Rows.Item(42).RowHeight = 45
$ws.Rows.Item(43).RowHeight = 52.2
$ws.Rows.Item(44).RowHeight = 47.4
$ws.Rows.Item(45).RowHeight = 45.6
$ws.Rows.Item(46).RowHeight = 51
$ws.Rows.Item(47).RowHeight = 47.4
$ws.Rows.Item(48).RowHeight = 46.2
$ws.Rows.Item(49).RowHeight = 57.6
$ws.Rows.Item(50).RowHeight = 50.4
$ws.Rows.Item(51).RowHeight = 58.8
$ws.Rows.Item(52).RowHeight = 57.6
$ws.Rows.Item(53).RowHeight = 51
$ws.Rows.Item(54).RowHeight = 49.8
$ws.Rows.Item(55).RowHeight = 49.2
$ws.Rows.Item(56).RowHeight = 48
$ws.Rows.Item(57).RowHeight = 46.8

# Header area vertical alignment: top -> center
$ws.Range("A1:G1").VerticalAlignment = -4108
$ws.Range("A2:G2").VerticalAlignment = -4108

# Restore selection near the edited region (best-effort; engine does not persist scroll/topLeftCell)
$ws.Range("D31").Select()
